# Update distractor analysis results (D column / B,C columns) per
# rerun with use_wle = TRUE (v.1.2.1), fixture ex3 distractors_summary.xlsx

$wb = $excel.ActiveWorkbook

# --- Sheet "correct" ---
$ws = $wb.Worksheets.Item("correct")
$updates = @{
    "D2" = 0.467
    "D3" = 0.548
    "D4" = 0.435
    "D5" = 0.458
    "D6" = 0.45
    "D7" = 0.531
    "D8" = 0.51
    "D9" = 0.537
    "D10" = 0.544
    "D11" = 0.54
    "D12" = 0.518
    "D13" = 0.511
    "D14" = 0.537
    "D15" = 0.514
    "D16" = 0.56
    "D17" = 0.521
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- Sheet "distractor" ---
$ws = $wb.Worksheets.Item("distractor")
$updates = @{
    "D2" = -0.24
    "D3" = -0.266
    "D4" = -0.246
    "D5" = -0.31
    "D6" = -0.199
    "D7" = -0.222
    "D8" = -0.167
    "D9" = -0.122
    "D10" = -0.092
    "D11" = -0.251
    "D12" = -0.253
    "D13" = -0.221
    "D14" = -0.207
    "D15" = -0.256
    "D16" = -0.233
    "D17" = -0.318
    "D18" = -0.239
    "D19" = -0.23
    "D20" = -0.236
    "D21" = -0.26
    "D22" = -0.252
    "D23" = -0.279
    "D24" = -0.251
    "D25" = -0.262
    "D26" = -0.248
    "D27" = -0.268
    "D28" = -0.244
    "D29" = -0.243
    "D30" = -0.211
    "D31" = -0.277
    "D32" = -0.24
    "D33" = -0.177
    "D34" = -0.262
    "D35" = -0.198
    "D36" = -0.204
    "D37" = -0.241
    "D38" = -0.177
    "D39" = -0.208
    "D40" = -0.268
    "D41" = -0.185
    "D42" = -0.183
    "D43" = -0.243
    "D44" = -0.133
    "D45" = -0.196
    "D46" = -0.248
    "D47" = -0.171
    "D48" = -0.172
    "D49" = -0.159
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- Sheet "descriptives" ---
$ws = $wb.Worksheets.Item("descriptives")
$updates = @{
    "B2" = 0.511
    "C2" = -0.224
    "B3" = 0.037
    "C3" = 0.046
    "B4" = 0.52
    "C4" = -0.24
    "B5" = 0.435
    "C5" = -0.318
    "B6" = 0.56
    "C6" = -0.092
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
